# Presença turma A. Script Cecon Covid
#
# Insert a new "presence" column M (header "2020-12-14"), shifting the
# existing "Email" column from M to N, and fill the new M column with the
# per-student numeric values below (three rows - the same ones whose
# "Resenha Regime de Metas" cell, L, was already blank - stay blank too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M; this shifts the old M (Email) to N and copies
# formatting (so M1 inherits the bold/border/centered header style, and the
# new M1 keeps style index "1" used by the other header cells).
$ws.Range("M1").EntireColumn.Insert()

# Write the new header text as a literal string ("2020-12-14") without
# letting Excel auto-convert it into a date serial value: build it via a
# text formula, then paste-special as values-only so no NumberFormat is
# ever touched (NumberFormat assignment would otherwise register a brand
# new style).
$ws.Range("M1").Formula = "=""2020-12-14"""
$ws.Range("M1").Copy()
$ws.Range("M1").PasteSpecial(-4163)  # xlPasteValues
$ws.Application.CutCopyMode = 0

# Per-student values for the new column M (row -> value). Rows not listed
# here (13, 37, 42) already had a blank "Resenha Regime de Metas" (L) cell
# in the source data and stay blank in the new column too - which is
# exactly the state EntireColumn.Insert() already leaves them in, so no
# further action is required for them.
$values = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 27.96
    7  = 1.08
    8  = 1.08
    9  = 0
    10 = 70.97
    11 = 2.15
    12 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 18.28
    18 = 0
    19 = 0
    20 = 36.56
    21 = 17.2
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 3.23
    28 = 0
    29 = 0
    30 = 0
    31 = 16.13
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    38 = 0
    39 = 0
    40 = 81.72
    41 = 0
    43 = 0
    44 = 10.75
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 13).Value = $values[$row]
}
